$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A (target stored OOXML width 16.42578125; engine quantizes
# ColumnWidth to 1/6-character steps, so 15.666666666666666 is the closest
# achievable setting, landing on a stored width of 16.5)
$ws.Columns.Item(1).ColumnWidth = 15.666666666666666

# New data values for A1:B32
$data = @(
    @(-0.064295737209334902,0.064008895235105001),
    @(0.035044934505112479,-0.035672774014473774),
    @(0.13860259853998258,-0.13922053783227639),
    @(-0.14476978134600316,0.14412070673814625),
    @(-0.13812070727812742,0.13683311301718692),
    @(-0.040457770617020827,0.040438006713475083),
    @(-0.020438007356913701,0.02041793884737686),
    @(-0.00041793949230051197,0.00040004899597789745),
    @(-0.038048316313604502,0.037850574482414245),
    @(-0.031850575046405538,0.031825063050696656),
    @(-0.027325063605704969,0.027281581871847038),
    @(-0.021281582436853963,0.021151925580523034),
    @(-0.039147128081051363,0.039082492959884974),
    @(-0.027082493569054122,0.027051345324177412),
    @(-0.021051345897775242,0.021026830601701185),
    @(-0.015026831176803146,0.015004283176844879),
    @(-0.0090042837541064458,0.0089999994037626152),
    @(-0.089053830305331871,0.088921956940129832),
    @(-0.0799219575004253,0.07887667828971523),
    @(-0.018013947379062145,0.018004354550580359),
    @(-0.0090043551276650646,0.0089999994224934099),
    @(-0.093929825083172602,0.093622704774215748),
    @(-0.084622705335425152,0.084124422048263092),
    @(-0.042124422822398699,0.041999999222049667),
    @(-0.090375342955940852,0.090208864979704373),
    @(-0.084208865537387823,0.083998062154488906),
    @(-0.077998062714602412,0.077291043003246962),
    @(-0.071291043574615465,0.070815172443280616),
    @(-0.058815173059286963,0.058596736743515621),
    @(-0.038596737412683435,0.038497533000890982),
    @(-0.023497533643647373,0.023483872611308243),
    @(-0.0060006365187499355,0.0059999994123254297)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
}
